$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.853.23'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.638.44'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.68'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5069'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2578'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06437'
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.58'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07770'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.287'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = '1.864.34'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').Value = '1.637.89'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5637'
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').Value = '0.0₅7616'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.14'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = '25.849.55'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.08'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.324'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.885'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.100'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.793'
$ws.Range('E25').Value = '  -5.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1275'
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '140.03'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.805'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.47'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.244'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04879'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.303'
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.223'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.560'
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.378'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9046'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.580'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('D38').Value = '1.131.10'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5515'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01563'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9960'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.535'
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.92'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('D45').Value = '1.774.74'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('E46').Value = '  -8.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.47'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4366'
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.722'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05049'
$ws.Range('E51').Value = '  +0.26%  '
